# Automatic update: refreshed MeteoCat daily-summary extraction timestamps
# and the handful of station readings that shifted between the 06:48-06:50
# run and the 07:18-07:20 run (2026-03-01 07:20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new literal text value.
# A leading apostrophe forces text entry for values that would otherwise
# be auto-converted to numbers by Excel (e.g. "98%" -> a percentage value);
# the apostrophe itself is not stored, matching the source data which is
# plain text in every cell of this sheet.
$updates = @(
    @{ Cell = 'E2'; Value = '2026-03-01 07:18:16' }
    @{ Cell = 'N2'; Value = '-2.5 °C 6:36 TU' }
    @{ Cell = 'O2'; Value = '-1.2 °C' }
    @{ Cell = 'E3'; Value = '2026-03-01 07:18:18' }
    @{ Cell = 'L3'; Value = '20.9 km/h - 131º 6:40 TU' }
    @{ Cell = 'N3'; Value = '-4.2 °C 6:48 TU' }
    @{ Cell = 'E4'; Value = '2026-03-01 07:18:21' }
    @{ Cell = 'H4'; Value = '''98%' }
    @{ Cell = 'E5'; Value = '2026-03-01 07:18:23' }
    @{ Cell = 'K5'; Value = '0.0 MJ/m2' }
    @{ Cell = 'N5'; Value = '-5.1 °C 6:36 TU' }
    @{ Cell = 'O5'; Value = '-4.0 °C' }
    @{ Cell = 'E6'; Value = '2026-03-01 07:18:26' }
    @{ Cell = 'H6'; Value = '''90%' }
    @{ Cell = 'E7'; Value = '2026-03-01 07:18:28' }
    @{ Cell = 'E8'; Value = '2026-03-01 07:18:30' }
    @{ Cell = 'E9'; Value = '2026-03-01 07:18:33' }
    @{ Cell = 'N9'; Value = '10.8 °C 6:59 TU' }
    @{ Cell = 'E10'; Value = '2026-03-01 07:18:35' }
    @{ Cell = 'K10'; Value = '0.0 MJ/m2' }
    @{ Cell = 'N10'; Value = '3.5 °C 6:47 TU' }
    @{ Cell = 'O10'; Value = '5.8 °C' }
    @{ Cell = 'E11'; Value = '2026-03-01 07:18:38' }
    @{ Cell = 'N11'; Value = '5.8 °C 6:54 TU' }
    @{ Cell = 'E12'; Value = '2026-03-01 07:18:40' }
    @{ Cell = 'H12'; Value = '''77%' }
    @{ Cell = 'O12'; Value = '10.1 °C' }
    @{ Cell = 'E13'; Value = '2026-03-01 07:18:42' }
    @{ Cell = 'H13'; Value = '''93%' }
    @{ Cell = 'N13'; Value = '2.8 °C 6:44 TU' }
    @{ Cell = 'O13'; Value = '4.1 °C' }
    @{ Cell = 'E14'; Value = '2026-03-01 07:18:44' }
    @{ Cell = 'E15'; Value = '2026-03-01 07:18:47' }
    @{ Cell = 'E16'; Value = '2026-03-01 07:18:49' }
    @{ Cell = 'H16'; Value = '''83%' }
    @{ Cell = 'K16'; Value = '0.0 MJ/m2' }
    @{ Cell = 'O16'; Value = '-5.4 °C' }
    @{ Cell = 'E17'; Value = '2026-03-01 07:18:52' }
    @{ Cell = 'E18'; Value = '2026-03-01 07:18:54' }
    @{ Cell = 'H18'; Value = '''100%' }
    @{ Cell = 'O18'; Value = '6.8 °C' }
    @{ Cell = 'E19'; Value = '2026-03-01 07:18:56' }
    @{ Cell = 'I19'; Value = '0.7 mm' }
    @{ Cell = 'N19'; Value = '5.8 °C 6:59 TU' }
    @{ Cell = 'E20'; Value = '2026-03-01 07:18:59' }
    @{ Cell = 'O20'; Value = '-3.3 °C' }
    @{ Cell = 'E21'; Value = '2026-03-01 07:19:01' }
    @{ Cell = 'J21'; Value = '1025.6 hPa' }
    @{ Cell = 'N21'; Value = '4.9 °C 6:34 TU' }
    @{ Cell = 'O21'; Value = '6.2 °C' }
    @{ Cell = 'E22'; Value = '2026-03-01 07:19:04' }
    @{ Cell = 'O22'; Value = '-5.5 °C' }
    @{ Cell = 'E23'; Value = '2026-03-01 07:19:06' }
    @{ Cell = 'N23'; Value = '-4.6 °C 6:59 TU' }
    @{ Cell = 'O23'; Value = '-3.7 °C' }
    @{ Cell = 'E24'; Value = '2026-03-01 07:19:08' }
    @{ Cell = 'O24'; Value = '4.8 °C' }
    @{ Cell = 'E25'; Value = '2026-03-01 07:19:11' }
    @{ Cell = 'H25'; Value = '''92%' }
    @{ Cell = 'E26'; Value = '2026-03-01 07:19:13' }
    @{ Cell = 'N26'; Value = '2.3 °C 6:47 TU' }
    @{ Cell = 'E27'; Value = '2026-03-01 07:19:16' }
    @{ Cell = 'G27'; Value = '159 cm' }
    @{ Cell = 'N27'; Value = '-2.9 °C 6:52 TU' }
    @{ Cell = 'E28'; Value = '2026-03-01 07:19:18' }
    @{ Cell = 'J28'; Value = '1025.6 hPa' }
    @{ Cell = 'N28'; Value = '8.2 °C 6:35 TU' }
    @{ Cell = 'E29'; Value = '2026-03-01 07:19:21' }
    @{ Cell = 'K29'; Value = '0.0 MJ/m2' }
    @{ Cell = 'N29'; Value = '5.9 °C 6:55 TU' }
    @{ Cell = 'O29'; Value = '8.8 °C' }
    @{ Cell = 'E30'; Value = '2026-03-01 07:19:23' }
    @{ Cell = 'H30'; Value = '''80%' }
    @{ Cell = 'E31'; Value = '2026-03-01 07:19:26' }
    @{ Cell = 'E32'; Value = '2026-03-01 07:19:28' }
    @{ Cell = 'L32'; Value = '7.2 km/h - 195º 6:38 TU' }
    @{ Cell = 'O32'; Value = '3.1 °C' }
    @{ Cell = 'E33'; Value = '2026-03-01 07:19:30' }
    @{ Cell = 'O33'; Value = '3.9 °C' }
    @{ Cell = 'E34'; Value = '2026-03-01 07:19:33' }
    @{ Cell = 'H34'; Value = '''99%' }
    @{ Cell = 'N34'; Value = '-0.8 °C 6:58 TU' }
    @{ Cell = 'O34'; Value = '-0.3 °C' }
    @{ Cell = 'E35'; Value = '2026-03-01 07:19:35' }
    @{ Cell = 'E36'; Value = '2026-03-01 07:19:38' }
    @{ Cell = 'H36'; Value = '''73%' }
    @{ Cell = 'M36'; Value = '12.8 °C 6:55 TU' }
    @{ Cell = 'O36'; Value = '10.4 °C' }
    @{ Cell = 'E37'; Value = '2026-03-01 07:19:40' }
    @{ Cell = 'N37'; Value = '6.0 °C 6:33 TU' }
    @{ Cell = 'O37'; Value = '6.2 °C' }
    @{ Cell = 'E38'; Value = '2026-03-01 07:19:42' }
    @{ Cell = 'E39'; Value = '2026-03-01 07:19:45' }
    @{ Cell = 'N39'; Value = '-4.5 °C 6:55 TU' }
    @{ Cell = 'O39'; Value = '-3.2 °C' }
    @{ Cell = 'E40'; Value = '2026-03-01 07:19:47' }
    @{ Cell = 'H40'; Value = '''90%' }
    @{ Cell = 'E41'; Value = '2026-03-01 07:19:49' }
    @{ Cell = 'O41'; Value = '11.5 °C' }
    @{ Cell = 'E42'; Value = '2026-03-01 07:19:52' }
    @{ Cell = 'H42'; Value = '''89%' }
    @{ Cell = 'N42'; Value = '5.2 °C 6:39 TU' }
    @{ Cell = 'O42'; Value = '7.9 °C' }
    @{ Cell = 'E43'; Value = '2026-03-01 07:19:54' }
    @{ Cell = 'N43'; Value = '8.2 °C 6:47 TU' }
    @{ Cell = 'E44'; Value = '2026-03-01 07:19:56' }
    @{ Cell = 'O44'; Value = '-3.0 °C' }
    @{ Cell = 'E45'; Value = '2026-03-01 07:19:59' }
    @{ Cell = 'N45'; Value = '0.6 °C 6:46 TU' }
    @{ Cell = 'O45'; Value = '3.1 °C' }
    @{ Cell = 'E46'; Value = '2026-03-01 07:20:01' }
    @{ Cell = 'J46'; Value = '1026.6 hPa' }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
